# Add a "different first page" footer setup to the document's (only)
# section:
#   - a blank default footer (applies to pages 2+)
#   - a first-page footer containing a right-aligned PAGE number field
#   - page numbering restarts at 0 on the title page
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Different first page -> turns on <w:titlePg/> and splits the footer
# story into a "default" slot and a "first" slot.
$sec.PageSetup.DifferentFirstPageHeaderFooter = $true

$footers = $sec.Footers

# wdHeaderFooterPrimary (1): default footer, left blank.
$fDefault = $footers.Item(1)
$fDefault.Range.InsertAfter("")

# wdHeaderFooterFirstPage (2): first-page footer gets a right-aligned
# PAGE field, numbering restarted at 0.
$fFirst = $footers.Item(2)
$pnums = $fFirst.PageNumbers
$pnums.Add(2)
$pnums.StartingNumber = 0
